# Issue 47280: LKSM: Trailing whitespace in Source name won't resolve when deriving samples
# Add trailing whitespace to the Name/StringCol values on rows 2 and 3 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "SampleSetBVT1  "
$ws.Range("C2").Value = "a "
$ws.Range("A3").Value = "SampleSetBVT2   "

$ws.Range("A2").Select()
